# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt"
# (Cilantro) immediately before the existing row 250. Excel's native row
# insert shifts row 250 (and everything below it, through row 294) down to
# row 251 (through row 295), which is exactly the effect shown in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(250).Insert()

$ws.Cells.Item(250, 1).Value = 4
$ws.Cells.Item(250, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(250, 3).Value = "Los Lagos"
$ws.Cells.Item(250, 4).Value = 44711
$ws.Cells.Item(250, 5).Value = 10
$ws.Cells.Item(250, 6).Value = 100112040
$ws.Cells.Item(250, 7).Value = "Cilantro"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 70
$ws.Cells.Item(250, 11).Value = 6000
$ws.Cells.Item(250, 12).Value = 6000
$ws.Cells.Item(250, 13).Value = 6000
$ws.Cells.Item(250, 14).Value = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(250, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(250, 16).Value = 3000
$ws.Cells.Item(250, 17).Value = 2
$ws.Cells.Item(250, 18).Value = "Hortaliza"
